$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "oxysoheraster euaster"
$ws.Range("A3").Value = "anthaster euaster"
$ws.Range("A1").Value = "style"
